# Auto-generated Excel COM-interop script to apply price-refresh updates
# to the Leviathan_Profits workbook (chore: update Sheets via scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 828.625
$ws.Range("I6").Value = 232.71428
$ws.Range("J6").Value = 5000
$ws.Range("K6").Value = 698.14284
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = -586.14284
$ws.Range("N6").Value = -15224
$ws.Range("H86").Value = 3749.75
$ws.Range("J86").Value = 3000
$ws.Range("L86").Value = 3000
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 3749.75
$ws.Range("J89").Value = 3000
$ws.Range("L89").Value = 15000
$ws.Range("N89").Value = -26232
$ws.Range("H92").Value = 500.69232
$ws.Range("I92").Value = 346.27274
$ws.Range("K92").Value = 346.27274
$ws.Range("M92").Value = 901.72726
$ws.Range("H100").Value = 4100
$ws.Range("J100").Value = 5000
$ws.Range("L100").Value = 5000
$ws.Range("N100").Value = -6082
$ws.Range("H132").Value = 1639.5405
$ws.Range("I132").Value = 1657.3055
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 4971.916499999999
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -2441.916499999999
$ws.Range("N132").Value = -8060
$ws.Range("H137").Value = 4315.273
$ws.Range("I137").Value = 3598.7715
$ws.Range("J137").Value = 5569.15
$ws.Range("K137").Value = 10796.3145
$ws.Range("L137").Value = 16707.45
$ws.Range("M137").Value = -8246.3145
$ws.Range("N137").Value = -21807.45

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 399.5
$ws.Range("I3").Value = 399.5
$ws.Range("K3").Value = 399.5
$ws.Range("M3").Value = -284.5
$ws.Range("H32").Value = 34502.43
$ws.Range("I32").Value = 7109.5454
$ws.Range("K32").Value = 7109.5454
$ws.Range("M32").Value = -6822.5454
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H44").Value = 15000
$ws.Range("J44").Value = 15000
$ws.Range("L44").Value = 15000
$ws.Range("N44").Value = -15976
$ws.Range("H45").Value = 350255.47
$ws.Range("I45").Value = 460691.12
$ws.Range("K45").Value = 460691.12
$ws.Range("M45").Value = -460314.12
$ws.Range("H61").Value = 1390
$ws.Range("I61").Value = 1030
$ws.Range("K61").Value = 1030
$ws.Range("M61").Value = -818
$ws.Range("H74").Value = 1910
$ws.Range("I74").Value = 1406.4706
$ws.Range("K74").Value = 1406.4706
$ws.Range("M74").Value = -532.4706000000001
$ws.Range("H77").Value = 1910
$ws.Range("I77").Value = 1406.4706
$ws.Range("K77").Value = 7032.353000000001
$ws.Range("M77").Value = -2664.353000000001
$ws.Range("H136").Value = 1390
$ws.Range("I136").Value = 1030
$ws.Range("K136").Value = 3090
$ws.Range("M136").Value = -540

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 34999.668
$ws.Range("J35").Value = 34999.668
$ws.Range("L35").Value = 34999.668
$ws.Range("N35").Value = -35619.668
$ws.Range("H105").Value = 2938.8
$ws.Range("J105").Value = 1399.6666
$ws.Range("L105").Value = 1399.6666
$ws.Range("N105").Value = -4893.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 127199.75
$ws.Range("J31").Value = 2624.5
$ws.Range("L31").Value = 2624.5
$ws.Range("N31").Value = -3214.5
$ws.Range("H34").Value = 127199.75
$ws.Range("J34").Value = 2624.5
$ws.Range("L34").Value = 2624.5
$ws.Range("N34").Value = -3028.5
$ws.Range("H99").Value = 15506.75
$ws.Range("J99").Value = 13788.4
$ws.Range("L99").Value = 13788.4
$ws.Range("N99").Value = -16784.4
$ws.Range("H126").Value = 15506.75
$ws.Range("J126").Value = 13788.4
$ws.Range("L126").Value = 41365.2
$ws.Range("N126").Value = -46305.2
$ws.Range("H132").Value = 4037.1936
$ws.Range("I132").Value = 3838.5217
$ws.Range("K132").Value = 11515.5651
$ws.Range("M132").Value = -8985.5651
$ws.Range("H134").Value = 2522.7805
$ws.Range("I134").Value = 2598.8
$ws.Range("K134").Value = 7796.400000000001
$ws.Range("M134").Value = -5261.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1985.4
$ws.Range("I3").Value = 2039.3334
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 6118.0002
$ws.Range("L3").Value = 4500
$ws.Range("M3").Value = -6006.0002
$ws.Range("N3").Value = -4724
$ws.Range("H45").Value = 11081.5
$ws.Range("J45").Value = 11081.5
$ws.Range("L45").Value = 33244.5
$ws.Range("N45").Value = -34308.5
$ws.Range("H68").Value = 1517.4
$ws.Range("I68").Value = 1517.4
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 4552.200000000001
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -3741.200000000001
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 1517.4
$ws.Range("I71").Value = 1517.4
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 13656.6
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -9600.6
$ws.Range("N71").ClearContents()
$ws.Range("H82").Value = 10505
$ws.Range("I82").Value = 8342
$ws.Range("K82").Value = 25026
$ws.Range("M82").Value = -24620
$ws.Range("H85").Value = 10505
$ws.Range("I85").Value = 8342
$ws.Range("K85").Value = 25026
$ws.Range("M85").Value = -23622
$ws.Range("H92").Value = 451.1
$ws.Range("I92").Value = 567.8
$ws.Range("J92").Value = 334.4
$ws.Range("K92").Value = 1703.4
$ws.Range("L92").Value = 1003.2
$ws.Range("M92").Value = -455.3999999999999
$ws.Range("N92").Value = -3499.2
$ws.Range("H140").Value = 3321.05
$ws.Range("I140").Value = 2839.7144
$ws.Range("K140").Value = 8519.143199999999
$ws.Range("M140").Value = -3339.143199999999
$ws.Range("H141").Value = 3142
$ws.Range("I141").Value = 3065.4546
$ws.Range("K141").Value = 9196.363799999999
$ws.Range("M141").Value = -4016.363799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5893333.5
$ws.Range("I11").Value = 9025000
$ws.Range("K11").Value = 9025000
$ws.Range("M11").Value = -9024861
$ws.Range("H12").Value = 1052.6316
$ws.Range("H70").Value = 4788
$ws.Range("I70").Value = 4748.364
$ws.Range("J70").Value = 4933.3335
$ws.Range("K70").Value = 4748.364
$ws.Range("L70").Value = 4933.3335
$ws.Range("M70").Value = -4478.364
$ws.Range("N70").Value = -5473.3335
$ws.Range("H73").Value = 4788
$ws.Range("I73").Value = 4748.364
$ws.Range("J73").Value = 4933.3335
$ws.Range("K73").Value = 4748.364
$ws.Range("L73").Value = 4933.3335
$ws.Range("M73").Value = -3812.364
$ws.Range("N73").Value = -6805.3335
$ws.Range("H132").Value = 1892.9375
$ws.Range("I132").Value = 2076.7693
$ws.Range("J132").Value = 1096.3334
$ws.Range("K132").Value = 6230.3079
$ws.Range("L132").Value = 3289.0002
$ws.Range("M132").Value = -3700.3079
$ws.Range("N132").Value = -8349.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2965.182
$ws.Range("I132").Value = 2488.0293
$ws.Range("K132").Value = 7464.0879
$ws.Range("M132").Value = -4934.0879
$ws.Range("H136").Value = 3432.9443
$ws.Range("I136").Value = 3008.08
$ws.Range("K136").Value = 9024.24
$ws.Range("M136").Value = -6474.24

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 4250
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 4250
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 4250
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -4480
$ws.Range("H8").Value = 1000
$ws.Range("I8").Value = 1000
$ws.Range("K8").Value = 1000
$ws.Range("M8").Value = -860
$ws.Range("H54").Value = 40999
$ws.Range("J54").Value = 49999
$ws.Range("L54").Value = 49999
$ws.Range("N54").Value = -51039
$ws.Range("H62").Value = 13259.5
$ws.Range("I62").Value = 11349.75
$ws.Range("J62").Value = 14532.667
$ws.Range("K62").Value = 11349.75
$ws.Range("L62").Value = 14532.667
$ws.Range("M62").Value = -10725.75
$ws.Range("N62").Value = -15780.667
$ws.Range("H65").Value = 13259.5
$ws.Range("I65").Value = 11349.75
$ws.Range("J65").Value = 14532.667
$ws.Range("K65").Value = 56748.75
$ws.Range("L65").Value = 72663.33499999999
$ws.Range("M65").Value = -53628.75
$ws.Range("N65").Value = -78903.33499999999
$ws.Range("H74").Value = 28738
$ws.Range("J74").Value = 28738
$ws.Range("L74").Value = 28738
$ws.Range("N74").Value = -30610
$ws.Range("H77").Value = 28738
$ws.Range("J77").Value = 28738
$ws.Range("L77").Value = 86214
$ws.Range("N77").Value = -95574
$ws.Range("H132").Value = 1365251
$ws.Range("I132").Value = 5836.577
$ws.Range("J132").Value = 5292448
$ws.Range("K132").Value = 17509.731
$ws.Range("L132").Value = 15877344
$ws.Range("M132").Value = -14979.731
$ws.Range("N132").Value = -15882404
$ws.Range("H136").Value = 1612.3334
$ws.Range("I136").Value = 894.1429000000001
$ws.Range("J136").Value = 2617.8
$ws.Range("K136").Value = 2682.4287
$ws.Range("L136").Value = 7853.400000000001
$ws.Range("M136").Value = -132.4287000000004
$ws.Range("N136").Value = -12953.4
